$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Posting Label")
Write-Host $ws.Name
$v = $ws.Cells.Item(2,2).Value
Write-Host ("B2=" + $v.ToString())
